# Tilføjede emner til tidsplan.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# F12: rename "Koderi af ..." -> "Kodning af ..."
$ws.Range("F12").Value = 'Kodning af "lavere klasser"'

# F13: new value with wrap-text style
$ws.Range("F13").Value = "Diskussion af controller-`nimplementation"

# F14: new value with wrap-text style
$ws.Range("F14").Value = "Implementation af `ntoString()-metoder i alle klasser"

# F15: new value with wrap-text style
$ws.Range("F15").Value = "Diskussion af GUI-`nmockups"

# Apply wrap text + alignment style to F13:F15 (matches new cellXfs with wrapText=1)
$wrapRange = $ws.Range("F13:F15")
$wrapRange.WrapText = $true
$wrapRange.HorizontalAlignment = -4131  # xlLeft
$wrapRange.VerticalAlignment = -4108    # xlCenter

# Row height adjustments (16.5 -> 16) for rows 3, 4, 7, 9, 10
$ws.Rows.Item(3).RowHeight = 16
$ws.Rows.Item(4).RowHeight = 16
$ws.Rows.Item(7).RowHeight = 16
$ws.Rows.Item(9).RowHeight = 16
$ws.Rows.Item(10).RowHeight = 16

# Update selection to match recorded state in the sheet view
$ws.Range("F16").Select() | Out-Null
